# Update the cryptocurrency price/volume table (Sheet1, columns B-E, rows 2-51).
# Values are plain text cells in the source data (coinranking.com scrape), so
# numeric-looking prices (e.g. "0.998") are written with a leading apostrophe to
# force Excel to keep them as text instead of auto-converting them to numbers -
# matching the original inlineStr/text cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.678.16'
$ws.Range("E2").Value = '  +1.35%  '

# Row 3
$ws.Range("D3").Value = '1.565.10'
$ws.Range("E3").Value = '  +0.11%  '

# Row 4
$ws.Range("D4").Value = '''0.998'
$ws.Range("E4").Value = '  -0.33%  '

# Row 5
$ws.Range("D5").Value = '''210.21'
$ws.Range("E5").Value = '  -0.12%  '

# Row 6
$ws.Range("D6").Value = '''0.511'
$ws.Range("E6").Value = '  +4.25%  '

# Row 7
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.31%  '

# Row 8
$ws.Range("D8").Value = '''24.99'
$ws.Range("E8").Value = '  +5.23%  '

# Row 9
$ws.Range("D9").Value = '''0.246'
$ws.Range("E9").Value = '  +0.83%  '

# Row 10
$ws.Range("E10").Value = '  +0.22%  '

# Row 11
$ws.Range("D11").Value = '''0.0895'
$ws.Range("E11").Value = '  +0.20%  '

# Row 12
$ws.Range("D12").Value = '1.788.95'
$ws.Range("E12").Value = '  +0.09%  '

# Row 13
$ws.Range("D13").Value = '1.562.82'
$ws.Range("E13").Value = '  -0.07%  '

# Row 14
$ws.Range("D14").Value = '28.710.55'
$ws.Range("E14").Value = '  +1.55%  '

# Row 15
$ws.Range("E15").Value = '  +1.16%  '

# Row 16
$ws.Range("E16").Value = '  -0.28%  '

# Row 17
$ws.Range("D17").Value = '''61.53'
$ws.Range("E17").Value = '  +0.91%  '

# Row 18
$ws.Range("D18").Value = '''230.06'
$ws.Range("E18").Value = '  +0.98%  '

# Row 19
$ws.Range("D19").Value = '''7.34'
$ws.Range("E19").Value = '  -0.07%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0683'
$ws.Range("E20").Value = '  +0.97%  '

# Row 21
$ws.Range("E21").Value = '  -0.47%  '

# Row 22
$ws.Range("D22").Value = '''3.93'
$ws.Range("E22").Value = '  +0.09%  '

# Row 23
$ws.Range("D23").Value = '''9.03'
$ws.Range("E23").Value = '  +1.58%  '

# Row 24
$ws.Range("D24").Value = '''2.07'
$ws.Range("E24").Value = '  +2.13%  '

# Row 25
$ws.Range("D25").Value = '''151.91'
$ws.Range("E25").Value = '  +1.16%  '

# Row 26
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").Value = '''0.106'
$ws.Range("E26").Value = '  +2.92%  '

# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '''14.79'
$ws.Range("E27").Value = '  -0.43%  '

# Row 28
$ws.Range("D28").Value = '''0.998'
$ws.Range("E28").Value = '  -0.33%  '

# Row 29
$ws.Range("E29").Value = '  -1.32%  '

# Row 30
$ws.Range("D30").Value = '''0.0459'
$ws.Range("E30").Value = '  -3.71%  '

# Row 31
$ws.Range("E31").Value = '  -1.93%  '

# Row 32
$ws.Range("E32").Value = '  +0.37%  '

# Row 33
$ws.Range("D33").Value = '1.401.38'
$ws.Range("E33").Value = '  +1.71%  '

# Row 34
$ws.Range("E34").Value = '  -3.05%  '

# Row 35
$ws.Range("E35").Value = '  -3.49%  '

# Row 36
$ws.Range("E36").Value = '  -1.75%  '

# Row 37
$ws.Range("D37").Value = '''2.68'
$ws.Range("E37").Value = '  +1.71%  '

# Row 38
$ws.Range("E38").Value = '  -2.06%  '

# Row 39
$ws.Range("D39").Value = '''0.0161'
$ws.Range("E39").Value = '  -0.52%  '

# Row 40
$ws.Range("E40").Value = '  +0.81%  '

# Row 41
$ws.Range("D41").Value = '''0.519'
$ws.Range("E41").Value = '  -0.14%  '

# Row 42
$ws.Range("D42").Value = '''0.998'
$ws.Range("E42").Value = '  -0.27%  '

# Row 43
$ws.Range("D43").Value = '''0.769'
$ws.Range("E43").Value = '  -1.45%  '

# Row 44
$ws.Range("E44").Value = '  -3.11%  '

# Row 45
$ws.Range("D45").Value = '''63.88'
$ws.Range("E45").Value = '  +2.80%  '

# Row 46
$ws.Range("D46").Value = '''5.23'
$ws.Range("E46").Value = '  -1.60%  '

# Row 47
$ws.Range("D47").Value = '1.702.27'
$ws.Range("E47").Value = '  +0.12%  '

# Row 48
$ws.Range("E48").Value = '  -4.98%  '

# Row 49
$ws.Range("D49").Value = '''84.83'
$ws.Range("E49").Value = '  -0.48%  '

# Row 50
$ws.Range("D50").Value = '''42.48'
$ws.Range("E50").Value = '  +4.97%  '

# Row 51
$ws.Range("D51").Value = '''0.0511'
$ws.Range("E51").Value = '  -0.64%  '
